$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "steering knuckle" deflections table (A13:E20) to K12:O19.
# This mirrors a cut + paste (drag/move) of the block ten columns to the
# right and one row up.
$src = $ws.Range("A13:E20")
$dst = $ws.Range("K12")
$src.Cut($dst)

# Clear out what remains of the old location (Cut leaves the emptied
# cells' styling behind) and drop the now-stale merged cell.
$ws.Range("A13:E13").UnMerge()
$ws.Range("A13:E20").Clear()

# Re-merge the header cell at its new home.
$ws.Range("K12:O12").Merge()

# Reflect the new selection left behind by the move.
$ws.Range("K12:O19").Select()
